$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly data rows (2-11) had their per-report values (Fecha, Volumen,
# Precio minimo/maximo/promedio, Origen, Precio $/Kg) reshuffled between
# rows, while the identifying columns (Mercado, Region, Categoria, etc.)
# stay identical across every row. Apply the new values directly per row.

$ws.Range("D2").Value = 44162
$ws.Range("J2").Value = 80
$ws.Range("K2").Value = 7000
$ws.Range("L2").Value = 8000
$ws.Range("M2").Value = 7562
$ws.Range("O2").Value = "Región de O'Higgins"
$ws.Range("P2").Value = 302

$ws.Range("D3").Value = 44160
$ws.Range("J3").Value = 80
$ws.Range("K3").Value = 6500
$ws.Range("L3").Value = 7000
$ws.Range("M3").Value = 6688
$ws.Range("O3").Value = "Región de O'Higgins"
$ws.Range("P3").Value = 268

$ws.Range("D4").Value = 44167
$ws.Range("J4").Value = 60
$ws.Range("K4").Value = 8000
$ws.Range("L4").Value = 9000
$ws.Range("M4").Value = 8500
$ws.Range("O4").Value = "Región del Maule"
$ws.Range("P4").Value = 340

$ws.Range("D5").Value = 44473
$ws.Range("J5").Value = 60
$ws.Range("K5").Value = 9500
$ws.Range("L5").Value = 10000
$ws.Range("M5").Value = 9750
$ws.Range("O5").Value = "Región del Maule"
$ws.Range("P5").Value = 390

$ws.Range("D6").Value = 44161
$ws.Range("J6").Value = 53
$ws.Range("K6").Value = 6500
$ws.Range("L6").Value = 7000
$ws.Range("M6").Value = 6764
$ws.Range("O6").Value = "Región de O'Higgins"
$ws.Range("P6").Value = 271

$ws.Range("D8").Value = 44448
$ws.Range("K8").Value = 14000
$ws.Range("L8").Value = 15000
$ws.Range("M8").Value = 14500
$ws.Range("O8").Value = "Provincia del Elquí"
$ws.Range("P8").Value = 580

$ws.Range("D9").Value = 44159
$ws.Range("J9").Value = 42
$ws.Range("K9").Value = 6500
$ws.Range("L9").Value = 7000
$ws.Range("M9").Value = 6738
$ws.Range("O9").Value = "Región del Maule"
$ws.Range("P9").Value = 270

$ws.Range("D10").Value = 44166
$ws.Range("J10").Value = 56
$ws.Range("K10").Value = 7500
$ws.Range("L10").Value = 8000
$ws.Range("M10").Value = 7804
$ws.Range("P10").Value = 312

$ws.Range("D11").Value = 44466
$ws.Range("J11").Value = 60
$ws.Range("K11").Value = 11000
$ws.Range("L11").Value = 12000
$ws.Range("M11").Value = 11500
$ws.Range("P11").Value = 460
